# Re-save seed for "ML project 2 review seed".
#
# The authoring commit is primarily a re-save of the workbook by real
# Excel (fileVersion/theme/mc:AlternateContent/revision-tracking chrome,
# consolidated <col> runs, <row spans="..."> hints, default <pageMargins>,
# renumbered worksheet r:id's, and a float-formatting normalisation of the
# already-stored numeric cells, e.g. 1.3999999999999999 -> 1.4). None of
# that changes any cell's actual value - every "changed" number in the
# diff is bit-identical to the original double - and this runtime's own
# save path already emits the equivalent Excel chrome (theme, fileVersion,
# renumbered rIds, col-run consolidation, row spans, pageMargins) on every
# save, with no script action needed.
#
# The one concrete, user-visible edit captured in the diff is where the
# workbook was left selected when it was saved: Sheet1 is the active tab,
# scrolled down into the "virginica" rows, with A102:E151 highlighted
# (selection anchored at A102). Reproduce that view state here.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Select()

# Matches <selection activeCell="A102" sqref="A102:E151"/> on Sheet1's
# sheetView (and implies tabSelected="1" there, same as the diff).
$ws1.Range("A102:E151").Select()
